# Add a new weekly price record for Feria Lagunitas de Puerto Montt - Cilantro.
# A new row is inserted at row 224 (pushing the existing rows 224:239 down to
# 225:240) and populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 224, shifting rows 224:239 -> 225:240.
$ws.Rows("224:224").Insert()

# Populate the newly inserted row 224 with the new record.
$ws.Cells.Item(224, 1).Value = 4
$ws.Cells.Item(224, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(224, 3).Value = 'Los Lagos'
$ws.Range("D224").Value = 44610
$ws.Range("D224").NumberFormat = $ws.Range("D225").NumberFormat
$ws.Cells.Item(224, 5).Value = 10
$ws.Cells.Item(224, 6).Value = 100112040
$ws.Cells.Item(224, 7).Value = 'Cilantro'
$ws.Cells.Item(224, 8).Value = 'Sin especificar'
$ws.Cells.Item(224, 9).Value = 'Primera'
$ws.Cells.Item(224, 10).Value = 220
$ws.Cells.Item(224, 11).Value = 16000
$ws.Cells.Item(224, 12).Value = 16000
$ws.Cells.Item(224, 13).Value = 16000
$ws.Cells.Item(224, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(224, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(224, 16).Value = 444
$ws.Cells.Item(224, 17).Value = 36
$ws.Cells.Item(224, 18).Value = 'Hortaliza'
